$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.681.95"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.918.36"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.31"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.44"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.23"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.136"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.70"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.91"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "3.379.48"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "2.926.20"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.976"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "51.707.80"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.57"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.82"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.58"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +12.26%  "
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").Value = "  +16.18%  "
$ws.Range("E30").Value = "  +14.66%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.14"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.06"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.34"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -15.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.46"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.79"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.76"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "2.133.39"
$ws.Range("E48").Value = "  -3.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.251"
$ws.Range("E49").Value = "  -7.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.10"
$ws.Range("E51").Value = "  -0.11%  "
